$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 25
$prevRow = 24

# Date column: force text so Excel doesn't auto-convert "2025/12/04" into a date serial.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025/12/04"
$ws.Cells.Item($row, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($row, 3).Value = 1327

# Match the formatting used by the rest of the data rows (centered alignment, default number format).
$newRange = $ws.Range("A" + $row + ":C" + $row)
$prevRange = $ws.Range("A" + $prevRow + ":C" + $prevRow)

$newRange.Style = "Normal"
$newRange.HorizontalAlignment = $prevRange.HorizontalAlignment
$newRange.VerticalAlignment = $prevRange.VerticalAlignment
